$wb = $excel.ActiveWorkbook

# Sheets "1er Parcial" and "3er Parcial" both receive the same updated
# statistics for rows 3 and 4 (Medina Tolentino Elio / 3ASV group).
$sheetNames = @("1er Parcial", "3er Parcial")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Row 3
    $ws.Range("E3").Value = 20
    $ws.Range("F3").Value = 13
    $ws.Range("G3").Value = 60.61
    $ws.Range("H3").Value = 39.39
    $ws.Range("I3").Value = 5.8
    $ws.Range("J3").Value = 0
    $ws.Range("K3").Value = 0

    # Row 4
    $ws.Range("E4").Value = 19
    $ws.Range("F4").Value = 14
    $ws.Range("G4").Value = 57.58
    $ws.Range("H4").Value = 42.42
    $ws.Range("I4").Value = 6.4
    $ws.Range("J4").Value = 0
    $ws.Range("K4").Value = 0
}
